# Update "想去人数" (interested-count) figures in column F for the two
# data sheets that carry the full listing ("展览" and "全部类型").
# "演出" and "本地生活" only contain the header row and are untouched.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F, as updated by the source site refresh.
$updates = @{
    2  = 628
    4  = 355
    7  = 114
    11 = 63
    12 = 115
    13 = 1110
    14 = 1457
    15 = 321
    16 = 359
    18 = 95
    21 = 96
    22 = 253
    23 = 277
    29 = 628
    30 = 301
    31 = 3936
    33 = 457
    34 = 230
    35 = 997
    36 = 97
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
